# "2nd Question type added along with survey deletion functionality"
#
# The DataSet sheet stores one survey definition per column: row 1 holds the
# field/question key, row 2 holds its value. Each question value now carries
# a "<questionType>|<questionText>[|option1|option2...]" prefix instead of
# being a bare question string, and a brand new 4th question (a checkbox
# question asking for gender) is appended in column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# Re-tag the existing questions with their question-type prefixes. G2/H2 are
# touched first so that F2 (whose text also gained Male/Female options) is
# the newest/last-appended shared string - matches the order the workbook
# was actually edited in.
$ws.Range("G2").Value = "singleTextbox|What is your favorite color?"
$ws.Range("H2").Value = "singleTextbox|What is the average speed of an unladen swallow?"
$ws.Range("F2").Value = "singleTextbox|What is your name?|Male|Female"

# New 4th survey question (a checkbox/multi-choice question).
$ws.Range("I1").Value = "question4"
$ws.Range("I2").Value = "checkbox|What is your gender?|Male|Female"

# Column widths were resized (best-fit) to accommodate the new/longer text.
$ws.Columns.Item(1).ColumnWidth = 8.666666666666666
$ws.Columns.Item(2).ColumnWidth = 29.5
$ws.Columns.Item(3).ColumnWidth = 21.833333333333332
$ws.Columns.Item(4).ColumnWidth = 10
$ws.Columns.Item(5).ColumnWidth = 26.333333333333332
$ws.Columns.Item(6).ColumnWidth = 45.166666666666664
$ws.Columns.Item(7).ColumnWidth = 38.833333333333336
$ws.Columns.Item(8).ColumnWidth = 59.833333333333336
$ws.Columns.Item(9).ColumnWidth = 42.333333333333336

# Scroll / selection position left over from editing the new column.
$ws.Range("J3").Select()
